$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.184.44"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.626.48"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.07"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.18"
$ws.Range("E6").Value = "  -2.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -4.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.633.38"
$ws.Range("E9").Value = "  +1.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.34"
$ws.Range("E10").Value = "  -5.32%  "

$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.342"
$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.079.53"
$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.187.79"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.26"
$ws.Range("E16").Value = "  -1.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.627.50"
$ws.Range("E18").Value = "  +0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.63"
$ws.Range("E19").Value = "  -2.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.26"
$ws.Range("E20").Value = "  -2.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.44"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.09"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.420"
$ws.Range("E25").Value = "  -1.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.988"
$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("E27").Value = "  -2.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0810"
$ws.Range("E28").Value = "  -3.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.06"
$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("E32").Value = "  -4.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.97"
$ws.Range("E33").Value = "  -2.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.51"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.95"
$ws.Range("E35").Value = "  -5.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.920"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("E37").Value = "  -4.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.861"
$ws.Range("E38").Value = "  +2.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.54"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("E40").Value = "  -3.91%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "291.94"
$ws.Range("E41").Value = "  +1.94%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("E42").Value = "  -3.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.627"
$ws.Range("E43").Value = "  +0.88%  "

$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0547"
$ws.Range("E46").Value = "  -2.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.42"
$ws.Range("E47").Value = "  -0.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.40"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.67"
$ws.Range("E50").Value = "  -3.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.962.04"
$ws.Range("E51").Value = "  +0.11%  "
